$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new helpline log entries as rows 22 and 23

$ws.Cells.Item(22,1).Value = "2025-06-16T18:03:53.816246"
$ws.Cells.Item(22,2).Value = "Abhinab"
$ws.Cells.Item(22,3).Value = "8876BX"
$ws.Cells.Item(22,4).Value = ""
$ws.Cells.Item(22,5).Value = ""
$ws.Cells.Item(22,6).Value = "portal nhi khul rha"
$ws.Cells.Item(22,7).Value = "Patna"

$ws.Cells.Item(23,1).Value = "2025-06-16T20:27:06.038859"
$ws.Cells.Item(23,2).Value = "Ashu Prasad"
$ws.Cells.Item(23,3).Value = "IONCBX"
$ws.Cells.Item(23,4).Value = ""
$ws.Cells.Item(23,5).Value = ""
$ws.Cells.Item(23,6).Value = "unable to get payment"
$ws.Cells.Item(23,7).Value = "Patna"
